$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns retain their exact string representation
# (avoids Excel auto-coercing numeric-looking strings like "1.00" -> 1)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.953.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.055.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.445"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.36%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.572.26"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +14.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.897.18"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.048.31"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.73%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.499"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.175"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0965"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.80%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.33"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0694"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.085.63"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.71"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.347.38"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.42%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.660"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.16%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0248"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.15"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0894"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.33%  "
